$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the Betarraga data block (row 214),
# pushing the existing rows 214-271 down to 216-273.
$ws.Rows("214:215").Insert()

# New row 214: "Primera" quality entry for the newly-added week (44663).
$ws.Cells.Item(214, 1).Value = 8
$ws.Cells.Item(214, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(214, 3).Value = "Coquimbo"
$ws.Cells.Item(214, 4).Value = 44663
$ws.Cells.Item(214, 5).Value = 4
$ws.Cells.Item(214, 6).Value = 100114014
$ws.Cells.Item(214, 7).Value = "Betarraga"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 2400
$ws.Cells.Item(214, 11).Value = 450
$ws.Cells.Item(214, 12).Value = 500
$ws.Cells.Item(214, 13).Value = 475
$ws.Cells.Item(214, 14).Value = '$/paquete 3 unidades'
$ws.Cells.Item(214, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(214, 16).Value = 158
$ws.Cells.Item(214, 17).Value = 3
$ws.Cells.Item(214, 18).Value = "Hortaliza"

# New row 215: "Segunda" quality entry for the same newly-added week (44663).
$ws.Cells.Item(215, 1).Value = 8
$ws.Cells.Item(215, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(215, 3).Value = "Coquimbo"
$ws.Cells.Item(215, 4).Value = 44663
$ws.Cells.Item(215, 5).Value = 4
$ws.Cells.Item(215, 6).Value = 100114014
$ws.Cells.Item(215, 7).Value = "Betarraga"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Segunda"
$ws.Cells.Item(215, 10).Value = 1560
$ws.Cells.Item(215, 11).Value = 350
$ws.Cells.Item(215, 12).Value = 400
$ws.Cells.Item(215, 13).Value = 375
$ws.Cells.Item(215, 14).Value = '$/paquete 3 unidades'
$ws.Cells.Item(215, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(215, 16).Value = 125
$ws.Cells.Item(215, 17).Value = 3
$ws.Cells.Item(215, 18).Value = "Hortaliza"
